# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn row gets its handback file/date populated, de-de row gets its handback
#    file/date populated (de-de had a slightly later handback time)
#  - "Latest Target File" column (I) gets a link to a.md on both language sheets
#  - A couple of columns get widened so the longer strings are readable

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a9423d9c9e9f522a22a1bf4e5011c45de41e37a/e2e/a.md"
$bUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a9423d9c9e9f522a22a1bf4e5011c45de41e37a/e2e/b.md"

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shows up on the Overview sheet as well as both language sheets)
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US", 1)
}

# ---------------------------------------------------------------------------
# 2. Fill in "Latest Target File" (I), "Latest Handback File" (J) and
#    "Latest Handback DateTime" (K) for rows 2 and 3 on both language sheets.
# ---------------------------------------------------------------------------

# zh-cn: handback finished 2016-09-03 16:43:47
$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Cells.Replace("0001-01-01 00:00:00", "2016-09-03 16:43:47", 1)

$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# de-de: handback finished a bit later, 2016-09-03 16:43:54
$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-03 16:43:54"

$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-03 16:43:54"

# ---------------------------------------------------------------------------
# 3. Hyperlinks: re-create A2/A3 (source file links) and add the new I2/I3
#    (target file) links, in the same left-to-right, top-to-bottom order so
#    relationship ids line up the way Excel would naturally allocate them.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $aUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $aUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $bUrl, "", "", "b.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $aUrl, "", "", "a.md")
}

# ---------------------------------------------------------------------------
# 4. Widen columns that now hold longer status/date/filename text.
#    (ColumnWidth is in character units; Excel quantizes to the nearest
#    pixel, so these values are chosen to land on the intended width.)
# ---------------------------------------------------------------------------
$wideStatusWidth = 29.166666666666668
$wideFileWidth   = 39.166666666666664

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth  # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth  # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth  = $wideStatusWidth   # C: Status
$wsZhCn.Columns.Item(10).ColumnWidth = $wideFileWidth     # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth  = $wideStatusWidth   # C: Status
$wsDeDe.Columns.Item(10).ColumnWidth = $wideFileWidth     # J: Latest Handback File
